# Commit: "unify the conception of DataNode, DataTable, Entity."
# The only content-level change in the workbook is renaming the single
# worksheet from "Property1" to "DataNode" (all other diff hunks are
# Excel-version/IDE metadata - fileVersion, xr/xr2/xr16 revision
# namespaces, absPath, uid stamps - that Excel stamps automatically on
# save and are not reachable/meaningful through the object model).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "DataNode"
